$wb = $excel.ActiveWorkbook

$oldVersion = "mines - version 1.0.0 (Feb 3 2026) (built on February 03 2026 10.14.00 EST)"
$newVersion = "Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on February 03 2026 17.29.55 EST)"

$aboutSheet = $wb.Worksheets.Item("About")
$boundariesSheet = $wb.Worksheets.Item("Boundaries and methane sources")

# A2: "Version: ..." line
$aboutSheet.Range("A2").Value = "Version: $newVersion"

# A6: Recommended Citation line
$aboutSheet.Range("A6").Value = "Recommended Citation:  ""Global Energy Monitor, Coal mine boundaries and methane sources for Uvalnaya Coal Mine, Russia, M1520, version '$newVersion'. (See the CC license for attribution requirements if sharing or adapting the data set.)"

# S2:S17 build_version column on the boundaries sheet
for ($row = 2; $row -le 17; $row++) {
    $boundariesSheet.Cells.Item($row, 19).Value = $newVersion
}
